$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.027.20'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.220.48'
$ws.Range("E3").Value = '  -1.24%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.628'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '74.12'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '44.30'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.20%  '
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.31%  '
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = '2.554.78'
$ws.Range("E14").Value = '  -1.13%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.847'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.64%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.13%  '
$ws.Range("D17").Value = '2.186.08'
$ws.Range("E17").Value = '  -3.08%  '
$ws.Range("D18").Value = '41.921.99'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("E19").Value = '  +11.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.15%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +31.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.46'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -6.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '11.60'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.38%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("E28").Value = '  -1.21%  '
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '166.69'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.65'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +13.61%  '
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  -4.38%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.11'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.24%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.23%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.03'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.33%  '
$ws.Range("B40").Value = 'MultiversX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '65.84'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.40%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.14'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.24%  '
$ws.Range("E43").Value = '  -1.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '104.38'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.94%  '
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("E47").Value = '  +6.12%  '
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.24%  '
$ws.Range("E50").Value = '  +0.60%  '
$ws.Range("D51").Value = '2.428.46'
$ws.Range("E51").Value = '  -1.20%  '
